$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target order of data rows (2-7): rotate so the last two rows move to the top
$data = @(
    @("operation_end_time", "datetime"),
    @("SubProcessID", "str"),
    @("stream:datastream", "dict"),
    @("time:timestamp", "datetime"),
    @("org:resource", "str"),
    @("concept:name", "str")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
